$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the sample number string (E7760 -> E7420). All of G2:G41 share
#    this one string, so updating the range updates the single shared-string entry.
$ws.Range("G2:G41").Value = "E7420"

# 2. Turn the literal boolean cells H2:H41 into real formulas "=FALSE()"
#    (was a stored boolean, now recomputed via a formula on every row).
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 8).Formula = "=FALSE()"
}

# 3. Restore the view: scroll back to the top and move the selection to the
#    G column instead of H.
$ws.Range("G2:G41").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
